{"js": "// Standardize the \"data\" placeholder to \"data_assinatura\" in the\n// \"{{ cidade_assinatura }}, {{ data }}.\" closing line of the procura\u00e7\u00e3o.\nconst body = context.document.body;\n\n// Locate the exact text \"{{ data }}\" (there is a single occurrence, right\n// after \"{{ cidade_assinatura }}, \").\nconst results = body.search(\"{{ data }}\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items/text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Could not find the text \"{{ data }}\" in the document body.');\n}\n\nconst target = results.items[0];\n// Replace just the found range's text, turning \"{{ data }}\" into\n// \"{{ data_assinatura }}\".\ntarget.insertText(\"{{ data_assinatura }}\", \"Replace\");\n\nawait context.sync();\n", "ps1": "# Standardize the \"data\" placeholder to \"data_assinatura\" in the\n# \"{{ cidade_assinatura }}, {{ data }}.\" closing line of the procura\u00e7\u00e3o.\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n$find.Text = \"{{ data }}\"\n$find.Replacement.Text = \"{{ data_assinatura }}\"\n\n# FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n# MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace\n$find.Execute($find.Text, $true, $false, $false, $false, $false, $true, \"wdFindContinue\", $false, $find.Replacement.Text, \"wdReplaceAll\")\n"}
